$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.678.63'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '3.164.46'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'529.00"
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = "'140.22"
$ws.Range("E6").Value = '  +1.36%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = '  +13.09%  '
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").Value = "'0.440"
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("E11").Value = '  +4.32%  '
$ws.Range("D12").Value = "'0.140"
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("D13").Value = '3.709.67'
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").Value = "'25.72"
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("E15").Value = '  +3.78%  '
$ws.Range("D16").Value = '58.716.45'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.163.82'
$ws.Range("E17").Value = '  +1.52%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'6.22"
$ws.Range("E18").Value = '  +2.87%  '
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = "'375.16"
$ws.Range("E21").Value = '  +3.58%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +4.62%  '
$ws.Range("D24").Value = "'69.73"
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D27").Value = "'8.44"
$ws.Range("E27").Value = '  +15.74%  '
$ws.Range("D28").Value = '0.0₃0861'
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = "'22.45"
$ws.Range("E29").Value = '  +5.10%  '
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("D34").Value = "'6.32"
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("D35").Value = "'156.78"
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("D37").Value = '2.706.33'
$ws.Range("E37").Value = '  +7.00%  '
$ws.Range("D38").Value = "'24.99"
$ws.Range("E38").Value = '  -1.98%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("D40").Value = "'0.0690"
$ws.Range("E40").Value = '  +3.08%  '
$ws.Range("E41").Value = '  +6.45%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = "'0.0293"
$ws.Range("E42").Value = '  +8.66%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = "'0.723"
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = "'39.17"
$ws.Range("E44").Value = '  +3.61%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").Value = '3.207.28'
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("E47").Value = '  +12.64%  '
$ws.Range("D48").Value = "'6.21"
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").Value = "'0.981"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("D50").Value = "'20.08"
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("E51").Value = '  +0.87%  '

# Reset style pointer on protected numeric-text cells so no extraneous style index is attached
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
